# Auto-generated edit script: updates Hyperion market-price figures
# across the per-job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 54902
$ws.Range("J93").Value = 54902
$ws.Range("L93").Value = 54902
$ws.Range("N93").Value = -59894
$ws.Range("H98").Value = 1370.75
$ws.Range("I98").Value = 1370.75
$ws.Range("K98").Value = 1370.75
$ws.Range("M98").Value = 127.25
$ws.Range("H100").Value = 300
$ws.Range("J100").Value = 300
$ws.Range("L100").Value = 300
$ws.Range("N100").Value = -1382
$ws.Range("H106").Value = 40004980
$ws.Range("I106").Value = 43483500
$ws.Range("K106").Value = 43483500
$ws.Range("M106").Value = -43482869
$ws.Range("H122").Value = 1370.75
$ws.Range("I122").Value = 1370.75
$ws.Range("K122").Value = 4112.25
$ws.Range("M122").Value = -1662.25
$ws.Range("H132").Value = 16397034
$ws.Range("I132").Value = 18521964
$ws.Range("K132").Value = 55565892
$ws.Range("M132").Value = -55563362
$ws.Range("H134").Value = 142163.72
$ws.Range("J134").Value = 142163.72
$ws.Range("L134").Value = 142163.72
$ws.Range("N134").Value = -152303.72
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = ""
$ws.Range("H138").Value = 3354.0876
$ws.Range("I138").Value = 1046.3636
$ws.Range("J138").Value = 4804.657
$ws.Range("K138").Value = 3139.0908
$ws.Range("L138").Value = 14413.971
$ws.Range("M138").Value = 2000.9092
$ws.Range("N138").Value = -24693.971
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = ""
$ws.Range("N139").Value = ""
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""
$ws.Range("H141").Value = 11581.333
$ws.Range("I141").Value = 5495.08
$ws.Range("J141").Value = 42012.6
$ws.Range("K141").Value = 16485.24
$ws.Range("L141").Value = 126037.8
$ws.Range("M141").Value = -11305.24
$ws.Range("N141").Value = -136397.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3536131.5
$ws.Range("I2").Value = 4714275
$ws.Range("J2").Value = 1700
$ws.Range("K2").Value = 4714275
$ws.Range("L2").Value = 1700
$ws.Range("M2").Value = -4714162
$ws.Range("N2").Value = -1926
$ws.Range("H116").Value = 3536131.5
$ws.Range("I116").Value = 4714275
$ws.Range("J116").Value = 1700
$ws.Range("K116").Value = 4714275
$ws.Range("L116").Value = 1700
$ws.Range("M116").Value = -4711981
$ws.Range("N116").Value = -6288
$ws.Range("H122").Value = 1738658.4
$ws.Range("I122").Value = 2498
$ws.Range("K122").Value = 7494
$ws.Range("M122").Value = -5044

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3536131.5
$ws.Range("I3").Value = 4714275
$ws.Range("J3").Value = 1700
$ws.Range("K3").Value = 4714275
$ws.Range("L3").Value = 1700
$ws.Range("M3").Value = -4714161
$ws.Range("N3").Value = -1928
$ws.Range("H20").Value = 2878.8125
$ws.Range("I20").Value = 2638.1667
$ws.Range("J20").Value = 3600.75
$ws.Range("K20").Value = 2638.1667
$ws.Range("L20").Value = 3600.75
$ws.Range("M20").Value = -2391.1667
$ws.Range("N20").Value = -4094.75
$ws.Range("H75").Value = 4998.25
$ws.Range("I75").Value = 3331
$ws.Range("K75").Value = 3331
$ws.Range("M75").Value = -2395
$ws.Range("H78").Value = 4998.25
$ws.Range("I78").Value = 3331
$ws.Range("K78").Value = 9993
$ws.Range("M78").Value = -5313
$ws.Range("H94").Value = 2278044.8
$ws.Range("I94").Value = 2564845
$ws.Range("K94").Value = 2564845
$ws.Range("M94").Value = -2564394
$ws.Range("H105").Value = 20840000
$ws.Range("I105").Value = 20840000
$ws.Range("K105").Value = 20840000
$ws.Range("M105").Value = -20838253

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 9005
$ws.Range("J21").Value = 9005
$ws.Range("L21").Value = 9005
$ws.Range("N21").Value = -9475
$ws.Range("H31").Value = 31239.334
$ws.Range("I31").Value = 8689.25
$ws.Range("K31").Value = 8689.25
$ws.Range("M31").Value = -8394.25
$ws.Range("H34").Value = 31239.334
$ws.Range("I34").Value = 8689.25
$ws.Range("K34").Value = 8689.25
$ws.Range("M34").Value = -8487.25
$ws.Range("H58").Value = 2145.2195
$ws.Range("I58").Value = 1619.129
$ws.Range("K58").Value = 1619.129
$ws.Range("M58").Value = -1416.129
$ws.Range("H134").Value = 43993.48
$ws.Range("I134").Value = 79403.38
$ws.Range("K134").Value = 238210.14
$ws.Range("M134").Value = -235675.14
$ws.Range("H135").Value = 148054.4
$ws.Range("J135").Value = 148054.4
$ws.Range("L135").Value = 148054.4
$ws.Range("N135").Value = -158194.4
$ws.Range("H136").Value = 2145.2195
$ws.Range("I136").Value = 1619.129
$ws.Range("K136").Value = 4857.387
$ws.Range("M136").Value = -2307.387

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 304
$ws.Range("I8").Value = 304
$ws.Range("K8").Value = 912
$ws.Range("M8").Value = -773
$ws.Range("H36").Value = 423
$ws.Range("I36").Value = 423
$ws.Range("K36").Value = 1269
$ws.Range("M36").Value = -1100
$ws.Range("H56").Value = 31255374
$ws.Range("I56").Value = 31255374
$ws.Range("K56").Value = 31255374
$ws.Range("M56").Value = -31254844
$ws.Range("H97").Value = 708
$ws.Range("J97").Value = 899.6667
$ws.Range("L97").Value = 2699.0001
$ws.Range("N97").Value = -3691.0001
$ws.Range("H132").Value = 2366.918
$ws.Range("I132").Value = 1294.9565
$ws.Range("J132").Value = 3015.7368
$ws.Range("K132").Value = 11654.6085
$ws.Range("L132").Value = 27141.6312
$ws.Range("M132").Value = -9124.6085
$ws.Range("N132").Value = -32201.6312

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 10042.071
$ws.Range("I132").Value = 7435
$ws.Range("J132").Value = 19601.334
$ws.Range("K132").Value = 22305
$ws.Range("L132").Value = 58804.00199999999
$ws.Range("M132").Value = -19775
$ws.Range("N132").Value = -63864.00199999999
$ws.Range("H134").Value = 41306.75
$ws.Range("J134").Value = 41306.75
$ws.Range("L134").Value = 123920.25
$ws.Range("N134").Value = -128990.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 48888.707
$ws.Range("I10").Value = 1724.7778
$ws.Range("K10").Value = 1724.7778
$ws.Range("M10").Value = -1584.7778
$ws.Range("H61").Value = 11121283
$ws.Range("I61").Value = 13898979
$ws.Range("K61").Value = 13898979
$ws.Range("M61").Value = -13898777
$ws.Range("H113").Value = 11121283
$ws.Range("I113").Value = 13898979
$ws.Range("K113").Value = 13898979
$ws.Range("M113").Value = -13896809
$ws.Range("H132").Value = 11894.41
$ws.Range("I132").Value = 12672.647
$ws.Range("J132").Value = 6602.4
$ws.Range("K132").Value = 38017.94100000001
$ws.Range("L132").Value = 19807.2
$ws.Range("M132").Value = -35487.94100000001
$ws.Range("N132").Value = -24867.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 8910.817999999999
$ws.Range("J31").Value = 8668.777
$ws.Range("L31").Value = 8668.777
$ws.Range("N31").Value = -9364.777
$ws.Range("H75").Value = 9000
$ws.Range("I75").Value = 9000
$ws.Range("K75").Value = 9000
$ws.Range("M75").Value = -8064
$ws.Range("H78").Value = 9000
$ws.Range("I78").Value = 9000
$ws.Range("K78").Value = 27000
$ws.Range("M78").Value = -22320
$ws.Range("H122").Value = 4597.8335
$ws.Range("I122").Value = 3397.5
$ws.Range("J122").Value = 6998.5
$ws.Range("K122").Value = 10192.5
$ws.Range("L122").Value = 20995.5
$ws.Range("M122").Value = -7742.5
$ws.Range("N122").Value = -25895.5
$ws.Range("H132").Value = 33711504
$ws.Range("I132").Value = 50016650
$ws.Range("J132").Value = 1101219.9
$ws.Range("K132").Value = 150049950
$ws.Range("L132").Value = 3303659.7
$ws.Range("M132").Value = -150047420
$ws.Range("N132").Value = -3308719.7
$ws.Range("H136").Value = 4431.884
$ws.Range("I136").Value = 4831.552
$ws.Range("K136").Value = 14494.656
$ws.Range("M136").Value = -11944.656

Write-Output "Updated $($wb.Worksheets.Count) worksheets."
